$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7838
$ws1.Range("F5").Value = 93
$ws1.Range("F6").Value = 209
$ws1.Range("F10").Value = 455
$ws1.Range("F16").Value = 26
$ws1.Range("F17").Value = 5717
$ws1.Range("F18").Value = 165
$ws1.Range("F19").Value = 238
$ws1.Range("F20").Value = 1371
$ws1.Range("F22").Value = 352

# Sheet "全部类型" (all types listing) - same events, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7838
$ws4.Range("F5").Value = 93
$ws4.Range("F6").Value = 209
$ws4.Range("F10").Value = 455
$ws4.Range("F16").Value = 26
$ws4.Range("F18").Value = 5717
$ws4.Range("F20").Value = 165
$ws4.Range("F21").Value = 238
$ws4.Range("F22").Value = 1371
$ws4.Range("F24").Value = 352
